# Aula 07 - Aplicando o componente de Card no nosso dashboard
#
# The original document doesn't have any w:proofErr markers (Word's
# spell-checker squiggly-range markers). The edit adds them around the
# words Word's PT-BR dictionary flags as misspelled (English/technical
# terms), which in turn forces those words into their own <w:r> runs
# (proofErr can only sit between runs). We rebuild each affected
# paragraph's XML directly via Range.InsertXML so we get the exact
# <w:proofErr>/<w:r> structure, and append the new "AULA 077" paragraphs
# the same way.

$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$indPpr = '<w:pPr><w:ind w:left="360"/></w:pPr>'

function Set-ParaXml($paraIndex, $innerXml) {
    $p = $d.Paragraphs.Item($paraIndex)
    $r = $p.Range
    $xml = "<w:p $wNs>$innerXml</w:p>"
    [void]$r.InsertXML($xml)
}

# 1) "Vuetify" - wrap whole word in spellStart/spellEnd
$inner1 = '<w:proofErr w:type="spellStart"/><w:r><w:t>Vuetify</w:t></w:r><w:proofErr w:type="spellEnd"/>'
Set-ParaXml 1 $inner1

# 3) "Elevation sobra de bordas" -> "Elevation" flagged, rest plain
$inner3 = $indPpr + '<w:proofErr w:type="spellStart"/><w:r><w:t>Elevation</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> sobra de bordas</w:t></w:r>'
Set-ParaXml 3 $inner3

# 4) "Display flex das div" -> "flex" and "div" flagged
$inner4 = $indPpr + '<w:r><w:t xml:space="preserve">Display </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>flex</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> das </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>div</w:t></w:r><w:proofErr w:type="spellEnd"/>'
Set-ParaXml 4 $inner4

# 10) "Cores e paletas de cores que o vuetify traz nativamente, segue a risca o material design."
#     -> "vuetify" and "a" flagged
$inner10 = $indPpr + '<w:r><w:t xml:space="preserve">Cores e paletas de cores que o </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>vuetify</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> traz nativamente, segue </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>a</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> risca o material design.</w:t></w:r>'
Set-ParaXml 10 $inner10

# 11) "Inversão de cores automativo." -> "automativo" flagged
$inner11 = $indPpr + '<w:r><w:t xml:space="preserve">Inversão de cores </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>automativo</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>.</w:t></w:r>'
Set-ParaXml 11 $inner11

# 14) "colunas" -> split (no proofErr) into "C" + "olunas"
$inner14 = $indPpr + '<w:r><w:t>C</w:t></w:r><w:r><w:t>olunas</w:t></w:r>'
Set-ParaXml 14 $inner14

# Append the new "AULA 077" paragraphs after the (now split) "colunas" paragraph.
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$insertPoint = $d.Range($lastPara.Range.End, $lastPara.Range.End)
$newPara1 = "<w:p $wNs>$indPpr</w:p>"
$newPara2 = "<w:p $wNs>$indPpr<w:r><w:t>AULA 077 – Card</w:t></w:r></w:p>"
$newPara3 = "<w:p $wNs>$indPpr<w:r><w:t>Cartões com imagem e textos informativos</w:t></w:r></w:p>"
$newParasXml = $newPara1 + $newPara2 + $newPara3
[void]$insertPoint.InsertXML($newParasXml)
